# Auto-update draw results: append the 2025-10-10 Pick 3 draw as a new row.
#
# The sheet is a flat results table (Date, Game, Phase, Result, InsertedAt)
# whose existing rows run from row 2 through row 23. This adds row 24 with
# the new day's draw, extending the table by exactly one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 24

# Columns A ("2025-10-10") and C ("251010") contain digit-only text that
# Excel would otherwise auto-convert to a date serial / plain number on
# assignment. Mark just those two cells as Text first so the values are
# stored verbatim as strings (matching the rest of the column), same as
# every prior row in this table.
$ws.Range("A$row").NumberFormat = "@"
$ws.Range("C$row").NumberFormat = "@"

$ws.Range("A$row").Value = "2025-10-10"
$ws.Range("B$row").Value = "Pick 3"
$ws.Range("C$row").Value = "251010"
$ws.Range("D$row").Value = "2-9-3"
$ws.Range("E$row").Value = "2025-10-10T21:37:00.349+04:00"
